$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = [double]"0.0002260325465254232"
$ws.Range("H2").Value = [double]"0.001205506914802257"
$ws.Range("K2").Value = [double]"5.104541448224953"
$ws.Range("L2").Value = "[1.7794937042566659, 8.42958919219324]"
$ws.Range("M2").Value = [double]"0.002724275993294167"
$ws.Range("N2").Value = [double]"0.002724275993294167"
$ws.Range("O2").Value = [double]"-0.5031579825569237"
$ws.Range("P2").Value = "[-1.1572633598809245, 0.1509473947670772]"
$ws.Range("Q2").Value = [double]"0.1311864383694625"
$ws.Range("R2").Value = [double]"0.1311864383694625"
$ws.Range("S2").Value = [double]"13.58401488854082"
$ws.Range("T2").Value = "[11.872715483617066, 15.29531429346457]"
$ws.Range("W2").Value = [double]"2.081281281281331"
$ws.Range("X2").Value = [double]"-0.624384384384399"
$ws.Range("Y2").Value = [double]"4.786946946947062"

# Row 3 updates
$ws.Range("E3").Value = [double]"23.65000000000026"
$ws.Range("G3").Value = [double]"2.065848421861638e-05"
$ws.Range("H3").Value = [double]"0.0004105286123198618"
$ws.Range("K3").Value = [double]"5.377381939508083"
$ws.Range("L3").Value = "[2.8536156819553655, 7.9011481970608]"
$ws.Range("M3").Value = [double]"3.418700811286257e-05"
$ws.Range("N3").Value = [double]"6.837401622572514e-05"
$ws.Range("O3").Value = [double]"1.591237119836271"
$ws.Range("P3").Value = "[0.9874475407679624, 2.19502669890458]"
$ws.Range("Q3").Value = [double]"3.421431047367918e-07"
$ws.Range("R3").Value = [double]"6.842862094735835e-07"
$ws.Range("S3").Value = [double]"13.83580393894077"
$ws.Range("T3").Value = "[12.242251235949084, 15.42935664193245]"
$ws.Range("W3").Value = [double]"17.66056056056075"
$ws.Range("X3").Value = [double]"15.38788788788806"
$ws.Range("Y3").Value = [double]"19.93323323323345"
